$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Build the value in a scratch cell forced to Text format so numeric-looking
    # strings (e.g. "43.07") are not auto-converted to numbers, then copy only the
    # value (not the scratch formatting) onto the destination cell.
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "D2" '66.946.62'
$ws.Range('E2').Value = '  +1.30%  '

Set-TextValue "D3" '3.896.86'
$ws.Range('E3').Value = '  +3.21%  '

Set-TextValue "D4" '1.00'
$ws.Range('E4').Value = '  -0.01%  '

Set-TextValue "D5" '471.97'
$ws.Range('E5').Value = '  +10.38%  '

Set-TextValue "D6" '143.60'
$ws.Range('E6').Value = '  +4.08%  '

$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('E9').Value = '  +0.60%  '

$ws.Range('E10').Value = '  +8.33%  '

$ws.Range('E11').Value = '  +10.33%  '

Set-TextValue "D12" '43.07'
$ws.Range('E12').Value = '  +1.65%  '

Set-TextValue "D13" '4.514.10'
$ws.Range('E13').Value = '  +2.91%  '

Set-TextValue "D14" '10.36'
$ws.Range('E14').Value = '  -0.10%  '

Set-TextValue "D15" '14.96'
$ws.Range('E15').Value = '  -0.29%  '

Set-TextValue "D16" '3.875.50'
$ws.Range('E16').Value = '  +2.91%  '

$ws.Range('E17').Value = '  -0.24%  '

$ws.Range('E18').Value = '  -0.03%  '

$ws.Range('E19').Value = '  +4.00%  '

Set-TextValue "D20" '67.154.94'
$ws.Range('E20').Value = '  +1.45%  '

Set-TextValue "D21" '430.61'
$ws.Range('E21').Value = '  +7.11%  '

Set-TextValue "D22" '3.38'
$ws.Range('E22').Value = '  +4.61%  '

Set-TextValue "D23" '14.62'
$ws.Range('E23').Value = '  -1.27%  '

$ws.Range('E24').Value = '  +4.64%  '

Set-TextValue "D25" '3.58'
$ws.Range('E25').Value = '  +9.65%  '

Set-TextValue "D26" '38.41'
$ws.Range('E26').Value = '  +5.12%  '

$ws.Range('E27').Value = '  +5.92%  '

Set-TextValue "D28" '10.03'
$ws.Range('E28').Value = '  +2.22%  '

Set-TextValue "D29" '9.61'
$ws.Range('E29').Value = '  -4.20%  '

Set-TextValue "D30" '727.27'
$ws.Range('E30').Value = '  +3.26%  '

Set-TextValue "D31" '13.69'
$ws.Range('E31').Value = '  +0.29%  '

$ws.Range('E32').Value = '  +0.32%  '

$ws.Range('E33').Value = '  +1.11%  '

Set-TextValue "D34" '43.28'
$ws.Range('E34').Value = '  +7.10%  '

Set-TextValue "D35" '0.155'
$ws.Range('E35').Value = '  +5.02%  '

Set-TextValue "D36" '57.15'
$ws.Range('E36').Value = '  +1.76%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('E38').Value = '  +20.47%  '

$ws.Range('E39').Value = '  -3.78%  '

Set-TextValue "D40" '0.0475'
$ws.Range('E40').Value = '  +1.19%  '

$ws.Range('E41').Value = '  +4.98%  '

$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('E43').Value = '  +4.69%  '

$ws.Range('E44').Value = '  +0.08%  '

$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D45" '2.55'
$ws.Range('E45').Value = '  -7.45%  '

$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D46" '2.78'
$ws.Range('E46').Value = '  +3.92%  '

$ws.Range('E47').Value = '  +6.41%  '

$ws.Range('E48').Value = '  +1.26%  '

Set-TextValue "D49" '3.19'
$ws.Range('E49').Value = '  -0.57%  '

Set-TextValue "D50" '144.14'
$ws.Range('E50').Value = '  +3.78%  '

$ws.Range('E51').Value = '  +4.21%  '

# Clean up the scratch cell used for forcing text values
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false
